$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gast"
$ws.Range("C2").Value = "Cckbr"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.806708
$ws.Range("H2").Value = 2.420124
$ws.Range("I2").Value = 0.3478655588626643
$ws.Range("J2").Value = 0.3478655588626643
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1226416666666667
$ws.Range("N2").Value = 0.367925
$ws.Range("O2").Value = 0.8154076983085706
$ws.Range("P2").Value = 0.8154076983085706
$ws.Range("Q2").Value = 0.09893601363333333
$ws.Range("R2").Value = 0.8904241226999999
$ws.Range("S2").Value = 0.2836522546730297
$ws.Range("T2").Value = 0.2836522546730297

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gast"
$ws.Range("C3").Value = "Cckbr"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.806708
$ws.Range("H3").Value = 2.420124
$ws.Range("I3").Value = 0.3478655588626643
$ws.Range("J3").Value = 0.3478655588626643
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.02776366666666667
$ws.Range("N3").Value = 0.083291
$ws.Range("O3").Value = 0.1845923016914294
$ws.Range("P3").Value = 0.1845923016914294
$ws.Range("Q3").Value = 0.02239717200933334
$ws.Range("R3").Value = 0.201574548084
$ws.Range("S3").Value = 0.06421330418963463
$ws.Range("T3").Value = 0.06421330418963463

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gast"
$ws.Range("C4").Value = "Cckbr"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.238656666666667
$ws.Range("H4").Value = 3.71597
$ws.Range("I4").Value = 0.534128821815285
$ws.Range("J4").Value = 0.5341288218152851
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1226416666666667
$ws.Range("N4").Value = 0.367925
$ws.Range("O4").Value = 0.8154076983085706
$ws.Range("P4").Value = 0.8154076983085706
$ws.Range("Q4").Value = 0.1519109180277778
$ws.Range("R4").Value = 1.36719826225
$ws.Range("S4").Value = 0.4355327531966702
$ws.Range("T4").Value = 0.4355327531966703

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gast"
$ws.Range("C5").Value = "Cckbr"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.238656666666667
$ws.Range("H5").Value = 3.71597
$ws.Range("I5").Value = 0.534128821815285
$ws.Range("J5").Value = 0.5341288218152851
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.02776366666666667
$ws.Range("N5").Value = 0.083291
$ws.Range("O5").Value = 0.1845923016914294
$ws.Range("P5").Value = 0.1845923016914294
$ws.Range("Q5").Value = 0.03438965080777778
$ws.Range("R5").Value = 0.30950685727
$ws.Range("S5").Value = 0.09859606861861482
$ws.Range("T5").Value = 0.09859606861861485

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Gast"
$ws.Range("C6").Value = "Cckbr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2736576666666666
$ws.Range("H6").Value = 0.820973
$ws.Range("I6").Value = 0.1180056193220505
$ws.Range("J6").Value = 0.1180056193220505
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1226416666666667
$ws.Range("N6").Value = 0.367925
$ws.Range("O6").Value = 0.8154076983085706
$ws.Range("P6").Value = 0.8154076983085706
$ws.Range("Q6").Value = 0.0335618323361111
$ws.Range("R6").Value = 0.302056491025
$ws.Range("S6").Value = 0.09622269043887056
$ws.Range("T6").Value = 0.09622269043887059

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Gast"
$ws.Range("C7").Value = "Cckbr"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2736576666666666
$ws.Range("H7").Value = 0.820973
$ws.Range("I7").Value = 0.1180056193220505
$ws.Range("J7").Value = 0.1180056193220505
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02776366666666667
$ws.Range("N7").Value = 0.083291
$ws.Range("O7").Value = 0.1845923016914294
$ws.Range("P7").Value = 0.1845923016914294
$ws.Range("Q7").Value = 0.007597740238111111
$ws.Range("R7").Value = 0.068379662143
$ws.Range("S7").Value = 0.02178292888317991
$ws.Range("T7").Value = 0.02178292888317992

